# Update finish_goods_stocks sheet: refresh stock snapshot (names reorder + new quantities)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product-name order for A2:A111 (same 110 products, re-ranked by current stock)
$names = @(
    'ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы',
    'Сб. Фитонефрол (Урологический сбор) 50г',
    'Эрва шерстистая трава 30г',
    'Ноготки цветки 50г',
    'Береза почки 50г',
    'Пустырник трава 50г',
    'Багульник болотный побеги 50г',
    'Чабрец трава 50г',
    'Боярышник плоды 75г',
    'Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г',
    'Валериана корневища с корнями 50г',
    'Шалфей листья 50г',
    'Подорожник большой листья 50г',
    'Рябина плоды 50г',
    'Бессмертник песчаный цветки 30г',
    'Девясил корневища и корни 50г',
    'Аир корневища 75г',
    'Укроп пахучий плоды 50г',
    'Сб. Грудной №4 50г',
    'Дуба кора 75г',
    'Лен семена 100г',
    'Крушина кора 50г',
    'Ламинарии слоевища (морская капуста) 100г',
    'Мать-и-мачеха листья 35г',
    'Ромашка цветки вн 50г',
    'Полынь горькая трава 50г',
    'Пижма цветки 75г',
    'Череда трава 50г',
    'Брусника листья 50г',
    'Шиповник плоды низковитаминные 50г',
    'Зверобой трава 50г',
    'Тысячелистник трава 50г',
    'Липа цветки 35г',
    'Кукуруза столбики с рыльцами 40г',
    'Эвкалипт прутовидный листья 75г',
    'Можжевельник плоды 50г',
    'Солодка корни 50г',
    'Толокнянка листья 50г',
    'Чага (березовый гриб) 50г',
    'Сенна листья 50г',
    'Спорыш трава 50г',
    'Крапива листья 50г',
    'Алтей корни 75г',
    'Сб. Фитопектол №1 (Грудной сбор №1) 35г',
    'Чистотел трава 50г',
    'Мята перечная листья 50г',
    'Сб. Фитопектол №2 (Грудной сбор №2) 35г',
    'Фп Детский травяной чай "ФармаЦветик® для иммунитета" 20х1,5 г',
    'Фп Фиточай "Лактафитол" (БАД) 20х1,5 г',
    'Фп Детский травяной чай "ФармаЦветик® для спокойного сна" 20х1,5 г',
    'Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г',
    'Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г',
    'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем"(БАД) 20*1,5г',
    'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем" (БАД) 20*1,5г',
    'Фп "Щедрость природы® Фиточай для иммунитета" 20х2,0 г',
    'Фп "Щедрость природы® Фиточай при простуде" 20х2,0 г',
    'Фп "Щедрость природы® Фиточай кардиологический" 20х2,0 г',
    'Фп "Щедрость природы® Фиточай успокоительный"20х2,0 г',
    'Фп "Щедрость природы® Фиточай диабетический" 20х2,0 г',
    'Фп Сб. Грудной №4 20x2,0г',
    'Фп Шалфей листья 20х1,5г',
    'Фп Сб. Бруснивер 20x2,0г',
    'Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г',
    'Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г',
    'Фп Чистотел трава 20х1,5г',
    'Фп Фиточай "Опалиховский" (БАД) 20х2,0 г',
    'Фп Фиточай "Тибетский" (БАД) 20х2,0  г',
    'Фп Мята перечная листья 20x1,5г',
    'Фп Сб. Арфазетин-Э 20x2,0г',
    'Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г',
    'Фп Сб. Элекасол 20x2,0г',
    'Фп Брусника листья 20х1,5г',
    'Фп Пустырник трава 20x1,5г',
    'Фп Подорожник листья 20x1,5г',
    'Фп Ромашка цветки 20x1,5г',
    'Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г',
    'Фп Мелисса лекарственная трава 20x1,5г',
    'Фп Сенна листья 20x1,5г',
    'Фп Пастушья сумка трава 20х1,5г',
    'Фп Череда трава 20х1,5г',
    'Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г',
    'Фп Шиповник плоды 20х2,0г',
    'Фп Зверобой трава 20x1,5г',
    'Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г',
    'Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г',
    'Фп Чабрец трава 20x1,5 г',
    'Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г',
    'Фп Душица трава 20x1,5г',
    'Фп Крапива листья 20x1,5г',
    'Фп Толокнянка листья 20x1,5г',
    'Фп Сб. Желудочный №3 20x2,0г',
    'Фп Хвощ полевой трава 20х1,5г',
    'Фп Липа цветки 20x1,5г',
    'Фп Береза листья 20x1,5г',
    'Фп Фиалка трехцветная трава 20x1,5г',
    'Фп Золототысячник трава 20х1,5г',
    'Фп Боярышник плоды 20х3,0г',
    'Фп Пижма цветки 20х1,5г',
    'Фп Аир корневища 20x1,5г',
    'Фп Ольха соплодия 20х1,5г',
    'Фп Лапчатка корневища 20x2,5г',
    'Фп Крушина кора 20x1,5г',
    'Фп Ноготки цветки 20x1,5г',
    'Фп Девясил корневища и корни 20х1,5г',
    'Фп Дуб кора 20х1,5г',
    'Фп Бадан корневища 20x1,5г',
    'Фп Валериана корневища с корнями 20x1,5г',
    'Фп Кровохлебка корневища и корни 20x1,5г',
    'Фп Тысячелистник трава 20x1,5г',
    'Фп Почечный чай листья 20x1,5г'
)

# New quantities for B2:B111 ($null marks the one now out-of-stock row)
$qtys = @(
    81376,
    3052,
    6871,
    14251,
    11309,
    7574,
    9301,
    13972,
    15674,
    3317,
    15916,
    29302,
    7381,
    1512,
    22907,
    14793,
    6121,
    55560,
    31444,
    61891,
    49828,
    9176,
    13870,
    24700,
    95294,
    38512,
    15654,
    11421,
    16853,
    36604,
    36904,
    16423,
    25304,
    31745,
    32345,
    15944,
    44913,
    9648,
    36512,
    28559,
    18797,
    19087,
    8255,
    7307,
    24738,
    35556,
    10088,
    $null,
    15121,
    3830,
    3840,
    5610,
    7700,
    9110,
    666,
    720,
    1242,
    1386,
    864,
    462852,
    112674,
    153157,
    67407,
    1512,
    26286,
    4554,
    9449,
    59464,
    35981,
    1638,
    38052,
    72251,
    39108,
    27170,
    1336014,
    23758,
    37962,
    73114,
    5164,
    49841,
    187033,
    51534,
    55157,
    91629,
    97429,
    75006,
    59941,
    31590,
    72659,
    44476,
    26365,
    31570,
    83499,
    6034,
    5152,
    5709,
    26018,
    11334,
    6503,
    5947,
    1813,
    5391,
    27959,
    11016,
    7245,
    1789,
    30932,
    10384,
    30884,
    157311
)

# Before touching B49/B102 (whose number formats swap), grab a format donor
# from a still-unmodified cell: B102 currently carries the "blank" style (s=5)
# and B2 carries the normal "#,##0" value style (s=2).
$ws.Range("B102").Copy()
$ws.Range("B49").PasteSpecial(-4122)   # xlPasteFormats -> B49 becomes the blank-cell style
$ws.Range("B2").Copy()
$ws.Range("B102").PasteSpecial(-4122)  # xlPasteFormats -> B102 becomes the normal value style
$excel.CutCopyMode = $false

for ($i = 0; $i -lt 110; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $qty = $qtys[$i]
    if ($null -eq $qty) {
        $ws.Cells.Item($row, 2).ClearContents()
    } else {
        $ws.Cells.Item($row, 2).Value = $qty
    }
}

# Restore the on-screen scroll/selection state recorded for this sheet.
$ws.Range("A94").Select()
try { $excel.ActiveWindow.ScrollRow = 79 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
